$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 31413.576
$ws.Cells.Item(98, 9).Value = 826.1875
$ws.Cells.Item(98, 10).Value = 80353.39999999999
$ws.Cells.Item(98, 11).Value = 826.1875
$ws.Cells.Item(98, 12).Value = 80353.39999999999
$ws.Cells.Item(98, 13).Value = 671.8125
$ws.Cells.Item(98, 14).Value = -83349.39999999999
$ws.Cells.Item(122, 8).Value = 31413.576
$ws.Cells.Item(122, 9).Value = 826.1875
$ws.Cells.Item(122, 10).Value = 80353.39999999999
$ws.Cells.Item(122, 11).Value = 2478.5625
$ws.Cells.Item(122, 12).Value = 241060.2
$ws.Cells.Item(122, 13).Value = -28.5625
$ws.Cells.Item(122, 14).Value = -245960.2
$ws.Cells.Item(126, 8).Value = 47006
$ws.Cells.Item(126, 10).Value = 47006
$ws.Cells.Item(126, 12).Value = 47006
$ws.Cells.Item(126, 14).Value = -56886
$ws.Cells.Item(138, 8).Value = 1601.99
$ws.Cells.Item(138, 9).Value = 973.69696
$ws.Cells.Item(138, 10).Value = 1911.4478
$ws.Cells.Item(138, 11).Value = 2921.09088
$ws.Cells.Item(138, 12).Value = 5734.3434
$ws.Cells.Item(138, 13).Value = 2218.90912
$ws.Cells.Item(138, 14).Value = -16014.3434
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 38994.11
$ws.Cells.Item(32, 9).Value = 45213.324
$ws.Cells.Item(32, 10).Value = 19771.092
$ws.Cells.Item(32, 11).Value = 45213.324
$ws.Cells.Item(32, 12).Value = 19771.092
$ws.Cells.Item(32, 13).Value = -44926.324
$ws.Cells.Item(32, 14).Value = -20345.092
$ws.Cells.Item(63, 8).Value = 2659.7917
$ws.Cells.Item(63, 9).Value = 2317.6316
$ws.Cells.Item(63, 10).Value = 3960
$ws.Cells.Item(63, 11).Value = 2317.6316
$ws.Cells.Item(63, 12).Value = 3960
$ws.Cells.Item(63, 13).Value = -1631.6316
$ws.Cells.Item(63, 14).Value = -5332
$ws.Cells.Item(66, 8).Value = 2659.7917
$ws.Cells.Item(66, 9).Value = 2317.6316
$ws.Cells.Item(66, 10).Value = 3960
$ws.Cells.Item(66, 11).Value = 11588.158
$ws.Cells.Item(66, 12).Value = 19800
$ws.Cells.Item(66, 13).Value = -8156.158000000001
$ws.Cells.Item(66, 14).Value = -26664
$ws.Cells.Item(102, 8).Value = 15528.134
$ws.Cells.Item(102, 9).Value = 1981.6666
$ws.Cells.Item(102, 11).Value = 1981.6666
$ws.Cells.Item(102, 13).Value = -359.6666
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).ClearContents()
$ws.Cells.Item(125, 14).Value = 0
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 544.3570999999999
$ws.Cells.Item(94, 9).Value = 427.36365
$ws.Cells.Item(94, 10).Value = 973.3333
$ws.Cells.Item(94, 11).Value = 427.36365
$ws.Cells.Item(94, 12).Value = 973.3333
$ws.Cells.Item(94, 13).Value = 23.63634999999999
$ws.Cells.Item(94, 14).Value = -1875.3333
$ws.Cells.Item(125, 8).Value = 50566
$ws.Cells.Item(125, 10).Value = 50566
$ws.Cells.Item(125, 12).Value = 50566
$ws.Cells.Item(125, 14).Value = -60406
$ws.Cells.Item(134, 8).Value = 2861.42
$ws.Cells.Item(134, 9).Value = 1074.3256
$ws.Cells.Item(134, 10).Value = 4209.579
$ws.Cells.Item(134, 11).Value = 3222.976799999999
$ws.Cells.Item(134, 12).Value = 12628.737
$ws.Cells.Item(134, 13).Value = -687.9767999999995
$ws.Cells.Item(134, 14).Value = -17698.737
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 981.7273
$ws.Cells.Item(16, 9).Value = 999.8333
$ws.Cells.Item(16, 11).Value = 999.8333
$ws.Cells.Item(16, 13).Value = -712.8333
$ws.Cells.Item(113, 8).Value = 981.7273
$ws.Cells.Item(113, 9).Value = 999.8333
$ws.Cells.Item(113, 11).Value = 999.8333
$ws.Cells.Item(113, 13).Value = 1170.1667
$ws.Cells.Item(115, 8).Value = 34256
$ws.Cells.Item(115, 10).Value = 34256
$ws.Cells.Item(115, 12).Value = 34256
$ws.Cells.Item(115, 14).Value = -36606
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(121, 8).Value = 282111.1
$ws.Cells.Item(121, 9).Value = 122
$ws.Cells.Item(121, 10).Value = 517102
$ws.Cells.Item(121, 11).Value = 366
$ws.Cells.Item(121, 12).Value = 1551306
$ws.Cells.Item(121, 13).Value = 944
$ws.Cells.Item(121, 14).Value = -1553926
$ws.Cells.Item(122, 8).Value = 5892.6
$ws.Cells.Item(122, 9).Value = 418.1111
$ws.Cells.Item(122, 11).Value = 3762.9999
$ws.Cells.Item(122, 13).Value = -1312.9999
$ws.Cells.Item(126, 8).Value = 26359.924
$ws.Cells.Item(126, 9).Value = 100486.664
$ws.Cells.Item(126, 10).Value = 4121.9
$ws.Cells.Item(126, 11).Value = 301459.992
$ws.Cells.Item(126, 12).Value = 12365.7
$ws.Cells.Item(126, 13).Value = -296519.992
$ws.Cells.Item(126, 14).Value = -22245.7
$ws.Cells.Item(131, 8).Value = 855.24
$ws.Cells.Item(131, 9).Value = 498.33334
$ws.Cells.Item(131, 10).Value = 903.9091
$ws.Cells.Item(131, 11).Value = 1495.00002
$ws.Cells.Item(131, 12).Value = 2711.7273
$ws.Cells.Item(131, 13).Value = 3544.99998
$ws.Cells.Item(131, 14).Value = -12791.7273
$ws.Cells.Item(137, 8).Value = 52643300
$ws.Cells.Item(137, 9).Value = 4226
$ws.Cells.Item(137, 10).Value = 71442970
$ws.Cells.Item(137, 11).Value = 12678
$ws.Cells.Item(137, 12).Value = 214328910
$ws.Cells.Item(137, 13).Value = -7578
$ws.Cells.Item(137, 14).Value = -214339110
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 929
$ws.Cells.Item(122, 9).Value = 853.0714
$ws.Cells.Item(122, 10).Value = 1141.6
$ws.Cells.Item(122, 11).Value = 2559.2142
$ws.Cells.Item(122, 12).Value = 3424.8
$ws.Cells.Item(122, 13).Value = -109.2142000000003
$ws.Cells.Item(122, 14).Value = -8324.799999999999
$ws.Cells.Item(132, 8).Value = 3705.303
$ws.Cells.Item(132, 9).Value = 1339.3529
$ws.Cells.Item(132, 10).Value = 6219.125
$ws.Cells.Item(132, 11).Value = 4018.0587
$ws.Cells.Item(132, 12).Value = 18657.375
$ws.Cells.Item(132, 13).Value = -1488.0587
$ws.Cells.Item(132, 14).Value = -23717.375
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1663.6666
$ws.Cells.Item(16, 9).Value = 2366.6667
$ws.Cells.Item(16, 10).Value = 960.6667
$ws.Cells.Item(16, 11).Value = 2366.6667
$ws.Cells.Item(16, 12).Value = 960.6667
$ws.Cells.Item(16, 13).Value = -2196.6667
$ws.Cells.Item(16, 14).Value = -1300.6667
$ws.Cells.Item(93, 8).Value = 1633.4762
$ws.Cells.Item(93, 9).Value = 1300.3
$ws.Cells.Item(93, 11).Value = 1300.3
$ws.Cells.Item(93, 13).Value = -52.29999999999995
$ws.Cells.Item(98, 8).Value = 48355
$ws.Cells.Item(98, 10).Value = 48355
$ws.Cells.Item(98, 12).Value = 48355
$ws.Cells.Item(98, 14).Value = -54345
$ws.Cells.Item(111, 8).Value = 43947.332
$ws.Cells.Item(111, 10).Value = 43947.332
$ws.Cells.Item(111, 12).Value = 43947.332
$ws.Cells.Item(111, 14).Value = -52127.332
$ws.Cells.Item(122, 8).Value = 38543.777
$ws.Cells.Item(122, 9).Value = 43149.25
$ws.Cells.Item(122, 10).Value = 1700
$ws.Cells.Item(122, 11).Value = 129447.75
$ws.Cells.Item(122, 12).Value = 5100
$ws.Cells.Item(122, 13).Value = -126997.75
$ws.Cells.Item(122, 14).Value = -10000
$ws.Cells.Item(132, 8).Value = 3847.6667
$ws.Cells.Item(132, 9).Value = 2249.5833
$ws.Cells.Item(132, 10).Value = 4557.926
$ws.Cells.Item(132, 11).Value = 6748.749899999999
$ws.Cells.Item(132, 12).Value = 13673.778
$ws.Cells.Item(132, 13).Value = -4218.749899999999
$ws.Cells.Item(132, 14).Value = -18733.778
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 498.33334
$ws.Cells.Item(100, 9).Value = 498.33334
$ws.Cells.Item(100, 11).Value = 996.66668
$ws.Cells.Item(100, 13).Value = -455.66668
$ws.Cells.Item(104, 8).Value = 40845.668
$ws.Cells.Item(104, 10).Value = 40845.668
$ws.Cells.Item(104, 12).Value = 40845.668
$ws.Cells.Item(104, 14).Value = -47833.668
$ws.Cells.Item(126, 8).Value = 2553.1
$ws.Cells.Item(126, 9).Value = 2371.4211
$ws.Cells.Item(126, 10).Value = 6005
$ws.Cells.Item(126, 11).Value = 7114.263300000001
$ws.Cells.Item(126, 12).Value = 18015
$ws.Cells.Item(126, 13).Value = -4644.263300000001
$ws.Cells.Item(126, 14).Value = -22955
$ws.Cells.Item(132, 8).Value = 1564.9131
$ws.Cells.Item(132, 9).Value = 1082.2812
$ws.Cells.Item(132, 11).Value = 3246.8436
$ws.Cells.Item(132, 13).Value = -716.8435999999997
$ws.Cells.Item(136, 8).Value = 18286.88
$ws.Cells.Item(136, 9).Value = 32014.719
$ws.Cells.Item(136, 11).Value = 96044.15700000001
$ws.Cells.Item(136, 13).Value = -93494.15700000001
